$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.737.27"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "1.852.78"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").Formula = "'0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Formula = "'263.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Formula = "'0.9990"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Formula = "'0.5396"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.54%  "
$ws.Range("D8").Formula = "'0.3211"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Formula = "'0.07090"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.13%  "
$ws.Range("D10").Formula = "'19.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").Formula = "'0.7822"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").Formula = "'0.07823"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "1.855.94"
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("D14").Formula = "'89.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Formula = "'5.063"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Formula = "'14.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("D17").Formula = "'0.9985"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Formula = "'0.000008038"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").Formula = "'0.9995"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "26.748.77"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "2.083.25"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").Formula = "'4.660"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Formula = "'6.070"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("D24").Formula = "'9.459"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").Formula = "'2.234"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("D26").Formula = "'142.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").Formula = "'1.701"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").Formula = "'17.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("D29").Formula = "'111.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Formula = "'4.305"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.06%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Formula = "'4.137"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Formula = "'0.08758"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Formula = "'0.04891"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").Formula = "'0.7406"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.73%  "
$ws.Range("D35").Formula = "'1.151"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").Formula = "'3.113"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("D38").Formula = "'2.378"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.44%  "
$ws.Range("D39").Formula = "'0.01760"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("D40").Formula = "'0.4866"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").Formula = "'0.9142"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").Formula = "'109.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D43").Formula = "'5.939"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").Formula = "'0.9988"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").Formula = "'7.757"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").Formula = "'0.4233"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Formula = "'9.138"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Formula = "'0.1260"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.13%  "
$ws.Range("D49").Formula = "'35.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").Formula = "'0.05846"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("D51").Formula = "'0.9014"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.58%  "
